$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K15").Value = 0.2970525035592049
$ws.Range("J16").Value = 0.2858677898194339
$ws.Range("I17").Value = 0.2775335613519331
$ws.Range("H18").Value = 0.2743085116504074
$ws.Range("G19").Value = 0.2534447081011285
$ws.Range("F20").Value = 0.2766837437271186
$ws.Range("E21").Value = 0.2867219094086165
$ws.Range("D22").Value = 0.1751453671933744
$ws.Range("C23").Value = 0.1965658720679752
$ws.Range("B24").Value = 0.4328090033804217
